$wb = $excel.ActiveWorkbook

# Sheet 1 (Zhanlan / Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 389
$ws.Range("F3").Value = 1077
$ws.Range("F4").Value = 9607
$ws.Range("F5").Value = 206
$ws.Range("F8").Value = 6557
$ws.Range("F10").Value = 10445
$ws.Range("F11").Value = 11587
$ws.Range("F13").Value = 1199
$ws.Range("F14").Value = 5030
$ws.Range("F15").Value = 830
$ws.Range("F16").Value = 491
$ws.Range("F17").Value = 98
$ws.Range("F20").Value = 1365
$ws.Range("F21").Value = 276
$ws.Range("F22").Value = 1903
$ws.Range("F23").Value = 920
$ws.Range("F24").Value = 1313
$ws.Range("F27").Value = 2073
$ws.Range("F28").Value = 444
$ws.Range("F29").Value = 660
$ws.Range("F30").Value = 2737
$ws.Range("F31").Value = 203
$ws.Range("F32").Value = 1825
$ws.Range("F34").Value = 828
$ws.Range("F35").Value = 84
$ws.Range("F36").Value = 935
$ws.Range("F37").Value = 21
$ws.Range("F38").Value = 53
$ws.Range("F39").Value = 3421
$ws.Range("F40").Value = 243
$ws.Range("F41").Value = 93
$ws.Range("F42").Value = 529
$ws.Range("F46").Value = 252
$ws.Range("F48").Value = 4236
$ws.Range("F49").Value = 80

# Sheet 2 (Yanchu / Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 13

# Sheet 3 (Bendi Shenghuo / Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 6075

# Sheet 4 (Quanbu Leixing / All types)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 389
$ws.Range("F3").Value = 1077
$ws.Range("F4").Value = 9608
$ws.Range("F6").Value = 13
$ws.Range("F9").Value = 10445
$ws.Range("F10").Value = 11587
$ws.Range("F12").Value = 1199
$ws.Range("F13").Value = 5030
$ws.Range("F14").Value = 830
$ws.Range("F15").Value = 491
$ws.Range("F16").Value = 98
$ws.Range("F18").Value = 0
$ws.Range("F20").Value = 1365
$ws.Range("F21").Value = 276
$ws.Range("F22").Value = 1903
$ws.Range("F23").Value = 920
$ws.Range("F24").Value = 1313
$ws.Range("F26").Value = 2073
$ws.Range("F27").Value = 444
$ws.Range("F28").Value = 660
$ws.Range("F29").Value = 2737
$ws.Range("F30").Value = 203
$ws.Range("F31").Value = 1825
$ws.Range("F34").Value = 828
$ws.Range("F38").Value = 84
$ws.Range("F39").Value = 935
$ws.Range("F40").Value = 21
$ws.Range("F42").Value = 243
$ws.Range("F43").Value = 93
$ws.Range("F44").Value = 529
$ws.Range("F47").Value = 252
$ws.Range("F49").Value = 4236
